# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The "Periodo Mora" (col E) / "Valor Mora" (col F) block on rows 16-29
# is reversed in place (newest period first), while every other column
# (B,C,D,G..J) and all row styles stay exactly where they are. Capture
# the current values first, then write them back in reverse order so
# each period keeps the "Valor Mora" it originally had.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 16
$lastRow = 29

$periodos = @()
$valores = @()

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $periodos += $ws.Cells.Item($r, 5).Value2
    $valores  += $ws.Cells.Item($r, 6).Value2
}

$count = $periodos.Count
for ($i = 0; $i -lt $count; $i++) {
    $targetRow = $firstRow + $i
    $sourceIndex = $count - 1 - $i
    $ws.Cells.Item($targetRow, 5).Value2 = $periodos[$sourceIndex]
    $ws.Cells.Item($targetRow, 6).Value2 = $valores[$sourceIndex]
}
